$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (shifts old rows 23-29 down to 24-30,
# carrying their existing content/styles along), giving us the new
# 30-row table shape.
$ws.Rows("23").Insert()

# --- Row 2: gauge changes from USACE/76065 to USGS/073802332 ---
$ws.Range("A2").Value = "USGS"
$ws.Range("B2").Value = "073802332"

# --- Row 7: gauge id changes from numeric 82742 to text "82740" ---
$ws.Range("B7").Value = "82740"

# --- Row 23 (new): USACE / 76305 ---
# B23 keeps the "Text" cell style (s=2 / numFmtId 49) after the row
# insert, which would coerce a plain .Value assignment into a shared
# string. Temporarily drop to the default style to store a real number,
# then restore the text style so the cell's appearance is unchanged.
$ws.Range("A23").Value = "USACE"
$ws.Range("B23").Style = "Normal"
$ws.Range("B23").Value = 76305
$ws.Range("B23").NumberFormat = "@"

# --- Column C (offset) updates ---
$ws.Range("C1").Value = 0.25
$ws.Range("C2").Value = 0.25
$ws.Range("C4").Value = 0.75
$ws.Range("C9").Value = 0.25
$ws.Range("C11").Value = 0.75
$ws.Range("C13").Value = 0.3
$ws.Range("C15").Value = 0.25
$ws.Range("C16").Value = 0.3
$ws.Range("C17").Value = 0.5
$ws.Range("C19").Value = 0.75
$ws.Range("C20").Value = 0.6
$ws.Range("C21").Value = 0.5
$ws.Range("C22").Value = 1.5
$ws.Range("C23").Value = 0.5
$ws.Range("C24").Value = -1.25
$ws.Range("C25").Value = -1.25
$ws.Range("C26").Value = -0.5
$ws.Range("C28").Value = -0.25
$ws.Range("C30").Value = -1.5

# --- Sheet view: scrolled to show row 8 onward, selection on C24 ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C24").Select()
